$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.960.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.07%  "

$ws.Range("D3").Value = "'1.638.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").Value = "'215.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("E8").Value = "  -0.71%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "'19.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("D11").Value = "'0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "'1.867.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").Value = "'1.623.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").Value = "'0.0₃0764"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "'62.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "'25.947.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("D19").Value = "'1.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").Value = "'192.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.40%  "

$ws.Range("E21").Value = "  -1.80%  "

$ws.Range("D22").Value = "'9.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.31%  "

$ws.Range("D23").Value = "'6.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.88%  "

$ws.Range("D24").Value = "'0.132"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.89%  "

$ws.Range("E25").Value = "  -0.41%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("D27").Value = "'143.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").Value = "'6.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").Value = "'15.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("D30").Value = "'1.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("D31").Value = "'0.0501"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("D32").Value = "'3.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E34").Value = "  -5.11%  "

$ws.Range("D35").Value = "'2.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.76%  "

$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("D37").Value = "'1.133.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("E38").Value = "  -2.00%  "

$ws.Range("D39").Value = "'2.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("E42").Value = "  -1.25%  "

$ws.Range("D43").Value = "'99.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.93%  "

$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").Value = "'1.777.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("E46").Value = "  +1.43%  "

$ws.Range("D47").Value = "'56.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("E48").Value = "  +2.17%  "

$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("D50").Value = "'7.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").Value = "'0.415"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
